$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, derived from the day's refreshed
# coinranking.com scrape (GitHub Actions cron job). All of these columns
# (Price / Volume(1h), plus the Coin name + Link for the two rows that
# swapped rank) are stored as plain text in the sheet, so every value is
# written back as a string and then forced to remain text (some look like
# plain numbers, e.g. "1.00", and Excel would otherwise silently coerce
# them to the Number type on assignment).
$updates = [ordered]@{
    "D2" = '26.318.51'
    "D3" = '1.585.03'
    "E3" = '  -1.03%  '
    "D4" = '1.00'
    "E4" = '  -0.04%  '
    "D5" = '209.57'
    "E5" = '  -0.94%  '
    "E6" = '  -1.51%  '
    "D7" = '1.00'
    "E7" = '  -0.04%  '
    "E8" = '  -1.09%  '
    "E9" = '  -0.35%  '
    "D10" = '19.53'
    "E10" = '  -1.03%  '
    "D11" = '0.0845'
    "E11" = '  +0.39%  '
    "D12" = '1.809.06'
    "D13" = '1.578.97'
    "E13" = '  -0.97%  '
    "E14" = '  -0.47%  '
    "D15" = '0.517'
    "E15" = '  -1.56%  '
    "D16" = '64.31'
    "E16" = '  -1.41%  '
    "D17" = '26.329.27'
    "E17" = '  -1.36%  '
    "E18" = '  -1.00%  '
    "D19" = '7.23'
    "E19" = '  -0.59%  '
    "D20" = '1.00'
    "E20" = '  -0.06%  '
    "D21" = '207.05'
    "E21" = '  -1.51%  '
    "E22" = '  -0.92%  '
    "E23" = '  -3.78%  '
    "D24" = '8.83'
    "E24" = '  -1.74%  '
    "D25" = '144.37'
    "E25" = '  +0.71%  '
    "D26" = '1.00'
    "E26" = '  -0.17%  '
    "D27" = '6.99'
    "E27" = '  -1.76%  '
    "E28" = '  -0.72%  '
    "D29" = '15.32'
    "E29" = '  -0.79%  '
    "D30" = '0.0503'
    "E30" = '  -2.13%  '
    "E31" = '  -0.73%  '
    "D32" = '3.24'
    "E32" = '  -1.24%  '
    "D33" = '2.95'
    "D34" = '1.31'
    "E34" = '  +14.26%  '
    "D35" = '1.282.79'
    "E35" = '  -1.20%  '
    "E36" = '  +0.46%  '
    "D37" = '0.606'
    "E37" = '  -1.10%  '
    "E38" = '  -1.30%  '
    "D39" = '0.0168'
    "E39" = '  -1.49%  '
    "E40" = '  -1.01%  '
    "E41" = '  +0.54%  '
    "E42" = '  -1.66%  '
    "E43" = '  -4.55%  '
    "D44" = '62.30'
    "E44" = '  -1.81%  '
    "D45" = '1.721.49'
    "E45" = '  -0.83%  '
    "D46" = '88.82'
    "E46" = '  -2.67%  '
    "E47" = '  -0.49%  '
    "E48" = '  +0.52%  '
    "E49" = '  -1.40%  '
    "B50" = 'BabyDogeCoin'
    "C50" = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    "D50" = '0.0₇0971'
    "E50" = '  -7.21%  '
    "B51" = 'USDD'
    "C51" = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
    "D51" = '1.00'
    "E51" = '  +0.01%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text interpretation so numeric-looking strings (e.g. "1.00",
    # "209.57") aren't auto-converted to the Number type by Excel, then
    # drop the temporary number-format override so the cell's style index
    # is left exactly as it was before (avoids introducing a new style).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}

Write-Output "Updated $($updates.Count) cells"
